$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.920.86'
$ws.Range("E2").Value = '  -2.05%  '
$ws.Range("D3").Value = '1.899.28'
$ws.Range("E3").Value = '  -4.00%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '324.36'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.4590'
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("D8").Value = '0.3808'
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("D9").Value = '0.07707'
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").Value = '  -4.21%  '
$ws.Range("D12").Value = '1.915.22'
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").Value = '6.916'
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("D14").Value = '5.642'
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").Value = '0.07050'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '83.76'
$ws.Range("E17").Value = '  -4.43%  '
$ws.Range("D18").Value = '0.000009469'
$ws.Range("E18").Value = '  -4.91%  '
$ws.Range("E19").Value = '  -3.95%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '28.888.17'
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("E22").Value = '  -5.15%  '
$ws.Range("D23").Value = '10.85'
$ws.Range("E23").Value = '  -3.01%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.150.32'
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '158.14'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '19.02'
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '5.616'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '117.31'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '1.839'
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.09244'
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.8558'
$ws.Range("E32").Value = '  -3.92%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.075'
$ws.Range("E33").Value = '  -2.96%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.234'
$ws.Range("E34").Value = '  -6.78%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.998'
$ws.Range("E35").Value = '  -6.16%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.05663'
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.138'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '1.004'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02027'
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5469'
$ws.Range("E40").Value = '  -4.62%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '7.377'
$ws.Range("E41").Value = '  -5.47%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1750'
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '9.258'
$ws.Range("E43").Value = '  -4.32%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '2.759'
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.5140'
$ws.Range("E45").Value = '  -4.20%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '11.12'
$ws.Range("E46").Value = '  -5.47%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.06812'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '2.057'
$ws.Range("E48").Value = '  -4.78%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.000002568'
$ws.Range("E49").Value = '  -16.56%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '110.01'
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.764'
$ws.Range("E51").Value = '  -3.57%  '
